$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 92, shifting the existing rows 92-128 down to 93-129.
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new weekly price record.
$ws.Cells.Item(92, 1).Value = 7
$ws.Cells.Item(92, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(92, 3).Value = 'Ñuble'
$ws.Cells.Item(92, 4).Value = 45029
$ws.Cells.Item(92, 5).Value = 16
$ws.Cells.Item(92, 6).Value = 100112037
$ws.Cells.Item(92, 7).Value = 'Cebollín'
$ws.Cells.Item(92, 8).Value = 'Sin especificar'
$ws.Cells.Item(92, 9).Value = 'Primera'
$ws.Cells.Item(92, 10).Value = 120
$ws.Cells.Item(92, 11).Value = 7000
$ws.Cells.Item(92, 12).Value = 7000
$ws.Cells.Item(92, 13).Value = 7000
$ws.Cells.Item(92, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(92, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(92, 16).Value = 194
$ws.Cells.Item(92, 17).Value = 36
$ws.Cells.Item(92, 18).Value = 'Hortaliza'
